# "make time the same" - correct two "days" values (column B) for patient K
# (rows 56-57) so that both measurements report the same underlying day value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B56").Value = 16.90123452
$ws.Range("B57").Value = 24.0575342465753

# Match the saved view/selection state: the active cell ends up on B57.
$ws.Range("B57").Select()
